$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels in D1 first, then C1, so the shared-string table
# ends up in the same order as the target workbook.
$ws.Range("D1").Value = "Composants"
$ws.Range("C1").Value = "Chef  Module"

# Widen columns C and D to match the target layout.
$ws.Columns.Item(3).ColumnWidth = 34.16666666666667
$ws.Columns.Item(4).ColumnWidth = 23.66666666666667

# Move the active selection to E8.
$ws.Range("E8").Select()
